$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Stock por almacén (separados por ; )"

# Column J needs to be widened to fit the new, longer header text.
$ws.Columns.Item(10).ColumnWidth = 32.85546875

# Update the view: scroll so column E is the left-most visible column,
# and select J3 instead of J4.
$ws.Range("J3").Select()
$excel.ActiveWindow.ScrollColumn = 5
